$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.838.86"
$ws.Range("E2").Value = "'  +0.14%  "
$ws.Range("D3").Value = "'1.742.08"
$ws.Range("E3").Value = "'  -0.95%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("D5").Value = "'225.26"
$ws.Range("E5").Value = "'  -5.12%  "
$ws.Range("D6").Value = "'0.9995"
$ws.Range("E6").Value = "'  -0.05%  "
$ws.Range("D7").Value = "'0.5146"
$ws.Range("E7").Value = "'  +1.45%  "
$ws.Range("D8").Value = "'0.2770"
$ws.Range("E8").Value = "'  +3.99%  "
$ws.Range("D9").Value = "'39.18"
$ws.Range("E9").Value = "'  -5.34%  "
$ws.Range("D10").Value = "'0.06095"
$ws.Range("E10").Value = "'  -1.71%  "
$ws.Range("D11").Value = "'1.734.28"
$ws.Range("E11").Value = "'  -1.36%  "
$ws.Range("D12").Value = "'0.06995"
$ws.Range("E12").Value = "'  +0.78%  "
$ws.Range("D13").Value = "'15.21"
$ws.Range("E13").Value = "'  -2.80%  "
$ws.Range("D14").Value = "'0.6331"
$ws.Range("E14").Value = "'  +4.51%  "
$ws.Range("D15").Value = "'4.504"
$ws.Range("E15").Value = "'  +0.67%  "
$ws.Range("D16").Value = "'76.56"
$ws.Range("E16").Value = "'  -1.30%  "
$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "'  -0.04%  "
$ws.Range("D18").Value = "'1.0000"
$ws.Range("E18").Value = "'  -0.01%  "
$ws.Range("D19").Value = "'25.852.40"
$ws.Range("E19").Value = "'  +0.05%  "
$ws.Range("D20").Value = "'11.45"
$ws.Range("E20").Value = "'  -1.45%  "
$ws.Range("D21").Value = "'0.000006625"
$ws.Range("E21").Value = "'  -2.95%  "
$ws.Range("D22").Value = "'1.958.23"
$ws.Range("E22").Value = "'  -1.12%  "
$ws.Range("D23").Value = "'4.087"
$ws.Range("E23").Value = "'  +0.37%  "
$ws.Range("D24").Value = "'8.491"
$ws.Range("E24").Value = "'  +3.76%  "
$ws.Range("D25").Value = "'5.103"
$ws.Range("E25").Value = "'  -1.81%  "
$ws.Range("D26").Value = "'137.35"
$ws.Range("E26").Value = "'  -0.50%  "
$ws.Range("D27").Value = "'1.502"
$ws.Range("E27").Value = "'  +3.01%  "
$ws.Range("D28").Value = "'1.815"
$ws.Range("E28").Value = "'  -0.36%  "
$ws.Range("D29").Value = "'14.98"
$ws.Range("E29").Value = "'  -0.34%  "
$ws.Range("D30").Value = "'102.74"
$ws.Range("E30").Value = "'  -0.02%  "
$ws.Range("D31").Value = "'0.08265"
$ws.Range("E31").Value = "'  +0.50%  "
$ws.Range("D32").Value = "'3.606"
$ws.Range("E32").Value = "'  -2.19%  "
$ws.Range("D33").Value = "'3.395"
$ws.Range("E33").Value = "'  -0.31%  "
$ws.Range("D34").Value = "'0.04404"
$ws.Range("E34").Value = "'  +0.70%  "
$ws.Range("D35").Value = "'2.618"
$ws.Range("E35").Value = "'  -1.39%  "
$ws.Range("D36").Value = "'0.9697"
$ws.Range("E36").Value = "'  -3.11%  "
$ws.Range("D37").Value = "'0.5968"
$ws.Range("E37").Value = "'  -1.50%  "
$ws.Range("D38").Value = "'2.647"
$ws.Range("E38").Value = "'  -3.07%  "
$ws.Range("E39").Value = "'  +0.17%  "
$ws.Range("D40").Value = "'1.921"
$ws.Range("E40").Value = "'  -1.02%  "
$ws.Range("D41").Value = "'0.9991"
$ws.Range("E41").Value = "'  -0.11%  "
$ws.Range("D42").Value = "'100.67"
$ws.Range("E42").Value = "'  -2.49%  "
$ws.Range("D43").Value = "'0.3818"
$ws.Range("E43").Value = "'  -0.22%  "
$ws.Range("D44").Value = "'0.7290"
$ws.Range("E44").Value = "'  -1.25%  "
$ws.Range("D45").Value = "'4.867"
$ws.Range("E45").Value = "'  -0.86%  "
$ws.Range("D46").Value = "'0.05472"
$ws.Range("E46").Value = "'  -0.40%  "
$ws.Range("D47").Value = "'6.250"
$ws.Range("D48").Value = "'0.1101"
$ws.Range("E48").Value = "'  +1.56%  "
$ws.Range("B49").Value = "'Aave"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'52.06"
$ws.Range("E49").Value = "'  -0.11%  "
$ws.Range("B50").Value = "'Elrond"
$ws.Range("C50").Value = "'https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'29.60"
$ws.Range("E50").Value = "'  -0.93%  "
$ws.Range("D51").Value = "'7.495"
$ws.Range("E51").Value = "'  -1.53%  "
